$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 71.72273850654022
$ws.Range("C2").Value = 75.92559237889883
$ws.Range("D2").Value = 66.96454883252153
$ws.Range("E2").Value = 81.26516350809368
$ws.Range("B3").Value = 94.79900830803332
$ws.Range("C3").Value = 92.94399963605767
$ws.Range("D3").Value = 94.23955928205713
$ws.Range("E3").Value = 94.198797909213
$ws.Range("B4").Value = 99.02226893752193
$ws.Range("C4").Value = 98.80816347945904
$ws.Range("D4").Value = 98.99522961424103
$ws.Range("E4").Value = 99.07492965207499
$ws.Range("B5").Value = 98.82380975909007
$ws.Range("C5").Value = 98.94682003618239
$ws.Range("D5").Value = 98.81328174316312
$ws.Range("E5").Value = 98.51446325242584
$ws.Range("B6").Value = 98.53829915692887
$ws.Range("C6").Value = 98.47172813241689
$ws.Range("D6").Value = 98.47431520005637
$ws.Range("E6").Value = 98.42970702669392
$ws.Range("B7").Value = 97.98219450712455
$ws.Range("C7").Value = 97.99486113053918
$ws.Range("D7").Value = 98.0413531971639
$ws.Range("E7").Value = 97.955660389679
$ws.Range("B8").Value = 97.5378865162842
$ws.Range("C8").Value = 97.49600345956333
$ws.Range("D8").Value = 97.55309504890855
$ws.Range("E8").Value = 97.50025393896046
$ws.Range("B9").Value = 96.16900119132742
$ws.Range("C9").Value = 96.16855075239933
$ws.Range("D9").Value = 96.15633004841125
$ws.Range("E9").Value = 96.20004453905575
